# Update "paises" (countries) workbook:
#  - refresh the "last updated" timestamp
#  - refresh a few existing countries' case counters
#  - insert Panama as a new row (pushing Sudafrica..Islandia down one row)
#  - replace the old Panama/Argentina row with refreshed Argentina data
#    (net effect: Argentina keeps its place near Serbia/Colombia, old
#    Panama row is gone, new Panama row appears higher up after Finlandia)
#  - swap Montenegro/Venezuela order and refresh Venezuela's counters

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 01:50"

# Estados Unidos (row 4): Casos totales, Nuevos casos, Casos activos, Recuperados
$ws.Range("B4").Value = 244230
$ws.Range("C4").Value = 29227
$ws.Range("D4").Value = 10403
$ws.Range("E4").Value = 227944

# Australia (row 23): Casos totales, Nuevos casos, Recuperados
$ws.Range("B23").Value = 5314
$ws.Range("C23").Value = 266
$ws.Range("E23").Value = 4701

# Rows 44-52 shift down to make room for the new Panama entry, and the
# trailing old Panama/Argentina rows collapse into a refreshed Argentina row.
$ws.Range("A44").Value = "Panama"
$ws.Range("B44").Value = 1475
$ws.Range("C44").Value = 158
$ws.Range("D44").Value = 9
$ws.Range("E44").Value = 1429
$ws.Range("F44").Value = 50
$ws.Range("G44").Value = 5
$ws.Range("H44").Value = 37

$ws.Range("A45").Value = "Sudafrica"
$ws.Range("B45").Value = 1462
$ws.Range("C45").Value = 82
$ws.Range("D45").Value = 50
$ws.Range("E45").Value = 1407
$ws.Range("F45").Value = 7
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 5

$ws.Range("A46").Value = "Peru"
$ws.Range("B46").Value = 1414
$ws.Range("C46").Value = 91
$ws.Range("D46").Value = 537
$ws.Range("E46").Value = 822
$ws.Range("F46").Value = 51
$ws.Range("G46").Value = 17
$ws.Range("H46").Value = 55

$ws.Range("A47").Value = "Republica Dominicana"
$ws.Range("B47").Value = 1380
$ws.Range("C47").Value = 96
$ws.Range("D47").Value = 16
$ws.Range("E47").Value = 1304
$ws.Range("F47").Value = 147
$ws.Range("G47").Value = 3
$ws.Range("H47").Value = 60

$ws.Range("A48").Value = "Mexico"
$ws.Range("B48").Value = 1378
$ws.Range("C48").Value = 163
$ws.Range("D48").Value = 35
$ws.Range("E48").Value = 1306
$ws.Range("F48").Value = 1
$ws.Range("G48").Value = 8
$ws.Range("H48").Value = 37

$ws.Range("A49").Value = "Islandia"
$ws.Range("B49").Value = 1319
$ws.Range("C49").Value = 99
$ws.Range("D49").Value = 284
$ws.Range("E49").Value = 1031
$ws.Range("F49").Value = 12
$ws.Range("G49").Value = 2
$ws.Range("H49").Value = 4

$ws.Range("A50").Value = "Argentina"
$ws.Range("B50").Value = 1265
$ws.Range("C50").Value = 132
$ws.Range("D50").Value = 256
$ws.Range("E50").Value = 973
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 4
$ws.Range("H50").Value = 36

$ws.Range("A51").Value = "Serbia"
$ws.Range("B51").Value = 1171
$ws.Range("C51").Value = 111
$ws.Range("D51").Value = 42
$ws.Range("E51").Value = 1098
$ws.Range("F51").Value = 81
$ws.Range("G51").Value = 3
$ws.Range("H51").Value = 31

$ws.Range("A52").Value = "Colombia"
$ws.Range("B52").Value = 1161
$ws.Range("C52").Value = 96
$ws.Range("D52").Value = 55
$ws.Range("E52").Value = 1087
$ws.Range("F52").Value = 50
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = 19

# Rows 109-110: Venezuela moves above Montenegro with refreshed counters;
# Montenegro keeps its old numbers but drops one row.
$ws.Range("A109").Value = "Venezuela"
$ws.Range("B109").Value = 146
$ws.Range("C109").Value = 2
$ws.Range("D109").Value = 43
$ws.Range("E109").Value = 98
$ws.Range("F109").Value = 6
$ws.Range("G109").Value = 2
$ws.Range("H109").Value = 5

$ws.Range("A110").Value = "Montenegro"
$ws.Range("B110").Value = 144
$ws.Range("C110").Value = 21
$ws.Range("D110").Value = 0
$ws.Range("E110").Value = 142
$ws.Range("F110").Value = 4
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 2
